$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Range("D3").Value = 6
$ws.Range("D4").Value = 0.5
